# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5 has a 3-column table whose table style is switched to a
#    different built-in style (GUID change only).
# 2) The deck's theme colour scheme (ppt/theme/theme2.xml, the theme used
#    by the slide master / whole deck) is swapped from the "Red Violet"
#    ("Integral") palette to the classic "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Update the table's style on slide 5 -------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{DBE7F589-F64F-4EC8-9751-2F15ACCD18BF}", $false)

# --- 2. Swap the theme colour scheme back to the "Office" palette --------
$themeSlide = $p.Slides.Item(1)
$tcs = $themeSlide.ThemeColorScheme

# Office theme colours, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# (values are VBA-style RGB() longs: R + G*256 + B*65536)
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
